$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New participant records to append (matches commit: 8 sample records total)
$data = @(
    @("P006", "Rina Sari",   "UNIV04", "rina.sari@univ4.edu",   "081678901234"),
    @("P007", "Andi Wijaya", "UNIV05", "andi.wijaya@univ5.edu", "081789012345"),
    @("P008", "Lisa Putri",  "UNIV06", "lisa.putri@univ6.edu",  "081890123456")
)

$startRow = 7
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    # Phone numbers carry a leading zero, so force the cell to Text
    # format before writing the value (otherwise Excel auto-converts
    # the entry to a number and the leading zero is lost).
    $ws.Cells.Item($row, 5).NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
